# repull data, push all data, mean calculation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the dSF (column F) values that changed after the data repull.
$ws.Range("F3").Value  = 6
$ws.Range("F5").Value  = 3
$ws.Range("F14").Value = -3
$ws.Range("F15").Value = 2
$ws.Range("F18").Value = 5
